$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(8, 17).Value = 1.67   # Q8: 1.69 -> 1.67
$ws.Cells.Item(9, 7).Value = 2.2   # G9: 2.22 -> 2.2
$ws.Cells.Item(9, 9).Value = 4.3   # I9: 4.4 -> 4.3
$ws.Cells.Item(9, 17).Value = 1.71   # Q9: 1.72 -> 1.71
$ws.Cells.Item(10, 6).Value = 5.2   # F10: 5.3 -> 5.2
$ws.Cells.Item(10, 7).Value = 5.3   # G10: 5.4 -> 5.3
$ws.Cells.Item(10, 14).Value = 4   # N10: 3.95 -> 4
$ws.Cells.Item(10, 15).Value = 1.3   # O10: 1.31 -> 1.3
$ws.Cells.Item(10, 16).Value = 2.04   # P10: 2.02 -> 2.04
$ws.Cells.Item(10, 17).Value = 1.93   # Q10: 1.94 -> 1.93
$ws.Cells.Item(10, 20).Value = 1.9   # T10: 1.92 -> 1.9
$ws.Cells.Item(10, 21).Value = 2.06   # U10: 2.04 -> 2.06
$ws.Cells.Item(10, 24).Value = 16   # X10: 15.5 -> 16
$ws.Cells.Item(10, 27).Value = 18.5   # AA10: 18 -> 18.5
$ws.Cells.Item(10, 28).Value = 18.5   # AB10: 18 -> 18.5
$ws.Cells.Item(10, 36).Value = 500   # AJ10: 140 -> 500
$ws.Cells.Item(10, 38).Value = 75   # AL10: 90 -> 75
$ws.Cells.Item(10, 41).Value = 10.5   # AO10: 11 -> 10.5
$ws.Cells.Item(11, 6).Value = 1.21   # F11: 1.97 -> 1.21
$ws.Cells.Item(11, 7).Value = 2.12   # G11: 2.64 -> 2.12
$ws.Cells.Item(11, 8).Value = 3.95   # H11: 3.4 -> 3.95
$ws.Cells.Item(11, 9).Value = 5.6   # I11: 5.2 -> 5.6
$ws.Cells.Item(11, 10).Value = 2.98   # J11: 2.72 -> 2.98
$ws.Cells.Item(11, 11).Value = 980   # K11: 5.1 -> 980
$ws.Cells.Item(11, 16).Value = 1.51   # P11: 1.42 -> 1.51
$ws.Cells.Item(11, 17).Value = 2.56   # Q11: 2.46 -> 2.56
$ws.Cells.Item(12, 6).Value = 2.22   # F12: 2.26 -> 2.22
$ws.Cells.Item(12, 8).Value = 3.75   # H12: 3.7 -> 3.75
$ws.Cells.Item(12, 11).Value = 3.35   # K12: 3.4 -> 3.35
$ws.Cells.Item(12, 16).Value = 1.7   # P12: 1.71 -> 1.7
$ws.Cells.Item(12, 17).Value = 2.02   # Q12: 2.22 -> 2.02
$ws.Cells.Item(13, 8).Value = 3.85   # H13: 3.8 -> 3.85
$ws.Cells.Item(13, 16).Value = 1.78   # P13: 1.77 -> 1.78
$ws.Cells.Item(14, 7).Value = 2.22   # G14: 2.24 -> 2.22
$ws.Cells.Item(14, 9).Value = 4.7   # I14: 4.8 -> 4.7
$ws.Cells.Item(16, 6).Value = 2.24   # F16: 2.2 -> 2.24
$ws.Cells.Item(16, 7).Value = 2.5   # G16: 2.68 -> 2.5
$ws.Cells.Item(16, 8).Value = 3.4   # H16: 3.45 -> 3.4
$ws.Cells.Item(16, 9).Value = 3.8   # I16: 4.4 -> 3.8
$ws.Cells.Item(16, 10).Value = 3.25   # J16: 2.9 -> 3.25
$ws.Cells.Item(16, 11).Value = 3.6   # K16: 3.65 -> 3.6
$ws.Cells.Item(16, 17).Value = 2.1   # Q16: 2.08 -> 2.1
